$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "317.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.95%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "47.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "10.67%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.285"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.17%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07933"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.42%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.595"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.09%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.324"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "32.53%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.642"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.36%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1274"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.36%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1934"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.85%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09433"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.24%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04632"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "10.95%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1045"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.13%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001333"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.65%"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.45%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005811"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.90%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.332"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.430"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.15%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3467"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.44%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.095"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.95%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1394"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.39%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3104"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.92%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001322"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.12%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004223"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-5.91%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001354"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.34%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003549"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-95.22%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02657"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "8.69%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05766"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "9.30%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01078"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "80.68%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008010"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.59%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1434"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.48%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007702"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.89%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008506"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "14.15%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3165"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.14%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006921"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.94%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.35%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05483"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "37.28%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004012"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.50%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.35%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.35%"
